$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 687/688, pushing the existing data (previously rows
# 687-734) down to rows 689-736. This mirrors the weekly data update where a
# new "Primera"/"Segunda" pair (fecha 44746) is added at the top of this
# block.
$ws.Range("A687:A688").EntireRow.Insert()

# Common (shared) values for both new rows.
$mercadoId = 8
$mercado = "Terminal La Palmera de La Serena"
$region = "Coquimbo"
$codreg = 4
$categoriaId = 100112008
$categoria = "Coliflor"
$variedad = "Sin especificar"
$unidad = "`$/unidad"
$origen = "Provincia del Elqu" + [char]237
$kgOUnidades = 1
$clasificacion = "Hortaliza"
$fecha = 44746

# Row 687: Calidad "Primera"
$ws.Cells.Item(687, 1).Value = $mercadoId
$ws.Cells.Item(687, 2).Value = $mercado
$ws.Cells.Item(687, 3).Value = $region
$ws.Cells.Item(687, 4).Value = $fecha
$ws.Cells.Item(687, 5).Value = $codreg
$ws.Cells.Item(687, 6).Value = $categoriaId
$ws.Cells.Item(687, 7).Value = $categoria
$ws.Cells.Item(687, 8).Value = $variedad
$ws.Cells.Item(687, 9).Value = "Primera"
$ws.Cells.Item(687, 10).Value = 2520
$ws.Cells.Item(687, 11).Value = 750
$ws.Cells.Item(687, 12).Value = 800
$ws.Cells.Item(687, 13).Value = 775
$ws.Cells.Item(687, 14).Value = $unidad
$ws.Cells.Item(687, 15).Value = $origen
$ws.Cells.Item(687, 16).Value = 775
$ws.Cells.Item(687, 17).Value = $kgOUnidades
$ws.Cells.Item(687, 18).Value = $clasificacion

# Row 688: Calidad "Segunda"
$ws.Cells.Item(688, 1).Value = $mercadoId
$ws.Cells.Item(688, 2).Value = $mercado
$ws.Cells.Item(688, 3).Value = $region
$ws.Cells.Item(688, 4).Value = $fecha
$ws.Cells.Item(688, 5).Value = $codreg
$ws.Cells.Item(688, 6).Value = $categoriaId
$ws.Cells.Item(688, 7).Value = $categoria
$ws.Cells.Item(688, 8).Value = $variedad
$ws.Cells.Item(688, 9).Value = "Segunda"
$ws.Cells.Item(688, 10).Value = 1400
$ws.Cells.Item(688, 11).Value = 650
$ws.Cells.Item(688, 12).Value = 700
$ws.Cells.Item(688, 13).Value = 675
$ws.Cells.Item(688, 14).Value = $unidad
$ws.Cells.Item(688, 15).Value = $origen
$ws.Cells.Item(688, 16).Value = 675
$ws.Cells.Item(688, 17).Value = $kgOUnidades
$ws.Cells.Item(688, 18).Value = $clasificacion
